$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33 (shifts old rows 33..71 down to 34..72).
$ws.Rows(33).Insert()

# New row 33 is a new weekly data point that sits (chronologically) between
# the old row 32 and old row 33 entries; it carries the same market/product
# metadata as row 32 did, but is a distinct record, while row 32's date
# moves forward to the newest sample (44601).
$ws.Range("A33").Value = 7
$ws.Range("B33").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C33").Value = "Ñuble"
$ws.Range("D33").Value = 44592
$ws.Range("E33").Value = 16
$ws.Range("F33").Value = "Fruta"
$ws.Range("G33").Value = 100108
$ws.Range("H33").Value = "Tropicales y subtropicales"
$ws.Range("I33").Value = 100108002
$ws.Range("J33").Value = "Mango"
$ws.Range("K33").Value = "Sin especificar"
$ws.Range("L33").Value = "Primera"
$ws.Range("M33").Value = 60
$ws.Range("N33").Value = 7000
$ws.Range("O33").Value = 7500
$ws.Range("P33").Value = 7250
$ws.Range("Q33").Value = "$/bandeja 4 kilos"
$ws.Range("R33").Value = "Perú"
$ws.Range("S33").Value = 1812
$ws.Range("T33").Value = 4

# Row 32 now reflects the newest weekly sample date.
$ws.Range("D32").Value = 44601
